# Update crypto price/volume figures per the Sun Jun 30 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.976.64'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = '  +1.57%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.423.25'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = '  +1.06%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.65'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  +1.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.02'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  +2.40%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  -0.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.59'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  -0.81%  '
$ws.Range("E10").Value = '  +0.88%  '
$ws.Range("E11").Value = '  -0.24%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.007.23'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E13").Value = '  -0.55%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.32'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  +1.84%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.418.26'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = '  +0.74%  '
$ws.Range("E16").Value = '  +0.00%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.006.95'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = '  +1.45%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.18'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = '  +1.26%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.00'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  +2.80%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '9.19'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  +3.45%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '390.60'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  +1.85%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '74.47'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  -1.21%  '
$ws.Range("E23").Value = '  +0.27%  '
$ws.Range("E24").Value = '  -0.12%  '
$ws.Range("E25").Value = '  +0.82%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.191'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  +3.37%  '
$ws.Range("E27").Value = '  +3.37%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  -0.08%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.05'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  +1.04%  '
$ws.Range("E30").Value = '  +0.67%  '
$ws.Range("E31").Value = '  +2.85%  '
$ws.Range("E32").Value = '  +0.02%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.58'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  +1.63%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.30'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  +6.56%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.98'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = '  +0.55%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '168.00'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  +1.24%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.455.37'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = '  +0.96%  '
$ws.Range("E38").Value = '  +0.67%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '28.61'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  +6.76%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0755'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  -1.42%  '
$ws.Range("E41").Value = '  +0.90%  '
$ws.Range("E42").Value = '  +1.99%  '
$ws.Range("E43").Value = '  +1.39%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.17'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  +4.69%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.527.82'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  +3.20%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.94'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  +0.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.63'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  -0.59%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.00'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  -0.02%  '
$ws.Range("E49").Value = '  +0.97%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.11'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  -1.45%  '
$ws.Range("E51").Value = '  -0.13%  '
